$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Foglio2")
$rng = $ws2.Range("L6")
$rng.Borders.Color = 0
$rng.Borders.Weight = 2
$rng.Borders.LineStyle = 1
$rng.Interior.PatternColorIndex = -4105
$rng.Interior.Pattern = -4142
